$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D26").Value = "다전공변경, 다전공을변경, 다전공바꾸, 다전공을바꾸, 다전공전환, 다전공을전환, 다전공교체, 다전공을교체, 다전공수정, 다전공을수정, 다전공변경하려"
$ws.Range("D27").Value = "복수전공변경, 복수전공을변경, 복수전공바꾸, 복수전공을바꾸, 복수전공전환, 복수전공을전환, 복수전공교체, 복수전공을교체, 복수전공수정, 복수전공을수정, 복수전공변경하려"
$ws.Range("D28").Value = "부전공변경, 부전공을변경, 부전공바꾸, 부전공을바꾸, 부전공전환, 부전공을전환, 부전공교체, 부전공을교체, 부전공수정, 부전공을수정, 부전공변경하려"
$ws.Range("D29").Value = "융합전공변경, 융합전공바꾸, 융합전공전환, 융합전공교체, 융합전공수정, 융합전공변경하려, 융합전공을변경, 융합전공을바꾸, 융합전공을전환, 융합전공을교체, 융합전공을수정"
$ws.Range("D30").Value = "융합부전공변경, 융합부전공바꾸, 융합부전공전환, 융합부전공교체, 융합부전공수정, 융합부전공변경하려, 융합부전공을변경, 융합부전공을바꾸, 융합부전공을전환, 융합부전공을교체, 융합부전공을수정"
$ws.Range("D31").Value = "연계전공변경, 연계전공바꾸, 연계전공전환, 연계전공교체, 연계전공수정, 연계전공변경하려, 연계전공을변경, 연계전공을바꾸, 연계전공을전환, 연계전공을교체, 연계전공을수정"
$ws.Range("D32").Value = "소단위전공과정변경, 소단위전공과정바꾸, 소단위전공과정전환, 소단위전공과정교체, 소단위전공과정수정, 소단위전공과정변경하려, 소단위전공과정을변경, 소단위전공과정을바꾸, 소단위전공과정을전환, 소단위전공과정을교체, 소단위전공과정을수정"
$ws.Range("D33").Value = "마이크로디그리변경, 마이크로디그리바꾸, 마이크로디그리전환, 마이크로디그리교체, 마이크로디그리수정, 마이크로디그리변경하려, 마이크로디그리을변경, 마이크로디그리을바꾸, 마이크로디그리을전환, 마이크로디그리을교체, 마이크로디그리을수정"
